$wb = $excel.ActiveWorkbook

# Delete the first "pageNavigation" sheet (sheetId=3), keep "pageNavigation (2)" (sheetId=4)
$excel.DisplayAlerts = $false
$wsOld = $wb.Worksheets.Item("pageNavigation")
$wsOld.Delete()
$excel.DisplayAlerts = $true

# Rename the remaining "pageNavigation (2)" sheet to "pageNavigation"
$ws = $wb.Worksheets.Item("pageNavigation (2)")
$ws.Name = "pageNavigation"

# Make it the active/selected sheet and set the new selection
$ws.Activate()
$ws.Range("C14").Select()
